# Translations check-in: rename several header labels on the "survey" and
# "settings" sheets to their new ".text" qualified forms, then leave the
# "settings" sheet as the active/selected tab (matching the new activeTab
# and tabSelected placement in the workbook).

$wb = $excel.ActiveWorkbook

# --- "survey" sheet header renames ---
$wsSurvey = $wb.Worksheets.Item("survey")
$wsSurvey.Range("J1").Value = "display.prompt.text"
$wsSurvey.Range("K1").Value = "display.hint.text"
$wsSurvey.Range("P1").Value = "display.constraint_message.text"
$wsSurvey.Range("S1").Value = "display.title.text"

# --- "settings" sheet header rename ---
$wsSettings = $wb.Worksheets.Item("settings")
$wsSettings.Range("C1").Value = "display.title.text"

# --- selection / active sheet state ---
$wsSurvey.Range("S2").Select()
$wsSettings.Activate()
$wsSettings.Range("C2").Select()
